$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4226.174
$ws.Range("I15").Value = 4226.174
$ws.Range("K15").Value = 12678.522
$ws.Range("M15").Value = -12509.522
$ws.Range("H116").Value = 2967.7104
$ws.Range("I116").Value = 2677.9167
$ws.Range("K116").Value = 2677.9167
$ws.Range("M116").Value = 764.0832999999998
$ws.Range("H137").Value = 1054.9143
$ws.Range("I137").Value = 911.9091
$ws.Range("J137").Value = 1296.9231
$ws.Range("K137").Value = 2735.7273
$ws.Range("L137").Value = 3890.7693
$ws.Range("M137").Value = -185.7273
$ws.Range("N137").Value = -8990.7693
$ws.Range("H138").Value = 1604.5424
$ws.Range("I138").Value = 1417.96
$ws.Range("J138").Value = 1741.7354
$ws.Range("K138").Value = 4253.88
$ws.Range("L138").Value = 5225.206200000001
$ws.Range("M138").Value = 886.1199999999999
$ws.Range("N138").Value = -15505.2062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 12725000
$ws.Range("I13").Value = 12725000
$ws.Range("K13").Value = 12725000
$ws.Range("M13").Value = -12724856
$ws.Range("H32").Value = 3062.6956
$ws.Range("I32").Value = 2778.3125
$ws.Range("K32").Value = 2778.3125
$ws.Range("M32").Value = -2491.3125
$ws.Range("H61").Value = 945.6667
$ws.Range("I61").Value = 945.6667
$ws.Range("K61").Value = 945.6667
$ws.Range("M61").Value = -733.6667
$ws.Range("H74").Value = 787.13043
$ws.Range("I74").Value = 460
$ws.Range("K74").Value = 460
$ws.Range("M74").Value = 414
$ws.Range("H77").Value = 787.13043
$ws.Range("I77").Value = 460
$ws.Range("K77").Value = 2300
$ws.Range("M77").Value = 2068
$ws.Range("H102").Value = 13891618
$ws.Range("I102").Value = 15154128
$ws.Range("K102").Value = 15154128
$ws.Range("M102").Value = -15152506
$ws.Range("H122").Value = 1776.5
$ws.Range("I122").Value = 1776.5
$ws.Range("K122").Value = 5329.5
$ws.Range("M122").Value = -2879.5
$ws.Range("H132").Value = 1576.7567
$ws.Range("I132").Value = 1265.6428
$ws.Range("K132").Value = 3796.9284
$ws.Range("M132").Value = -1266.9284
$ws.Range("H136").Value = 945.6667
$ws.Range("I136").Value = 945.6667
$ws.Range("K136").Value = 2837.0001
$ws.Range("M136").Value = -287.0001000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 683.1667
$ws.Range("I80").Value = 550.1667
$ws.Range("K80").Value = 550.1667
$ws.Range("M80").Value = 447.8333
$ws.Range("H83").Value = 683.1667
$ws.Range("I83").Value = 550.1667
$ws.Range("K83").Value = 2750.8335
$ws.Range("M83").Value = 2241.1665
$ws.Range("H99").Value = 71430240
$ws.Range("I99").Value = 83334950
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 83334950
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -83333452
$ws.Range("N99").Value = -4996
$ws.Range("H134").Value = 6196.3335
$ws.Range("I134").Value = 935.65
$ws.Range("K134").Value = 2806.95
$ws.Range("M134").Value = -271.9499999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1291.5
$ws.Range("I31").Value = 1024.4286
$ws.Range("J31").Value = 2092.7144
$ws.Range("K31").Value = 1024.4286
$ws.Range("L31").Value = 2092.7144
$ws.Range("M31").Value = -729.4286
$ws.Range("N31").Value = -2682.7144
$ws.Range("H34").Value = 1291.5
$ws.Range("I34").Value = 1024.4286
$ws.Range("J34").Value = 2092.7144
$ws.Range("K34").Value = 1024.4286
$ws.Range("L34").Value = 2092.7144
$ws.Range("M34").Value = -822.4286
$ws.Range("N34").Value = -2496.7144
$ws.Range("H58").Value = 1787.65
$ws.Range("I58").Value = 1450.2
$ws.Range("K58").Value = 1450.2
$ws.Range("M58").Value = -1247.2
$ws.Range("H107").Value = 543.9091
$ws.Range("I107").Value = 413.625
$ws.Range("K107").Value = 413.625
$ws.Range("M107").Value = 1506.375
$ws.Range("H136").Value = 1787.65
$ws.Range("I136").Value = 1450.2
$ws.Range("K136").Value = 4350.6
$ws.Range("M136").Value = -1800.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14286924
$ws.Range("J131").Value = 1266.8485
$ws.Range("L131").Value = 3800.5455
$ws.Range("N131").Value = -13880.5455
$ws.Range("H139").Value = 2529.889
$ws.Range("I139").Value = 2349.077
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 7047.231000000001
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -1907.231000000001
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 2050
$ws.Range("I13").Value = 2000
$ws.Range("J13").Value = 2100
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 2100
$ws.Range("M13").Value = -1861
$ws.Range("N13").Value = -2378
$ws.Range("H122").Value = 1482.6666
$ws.Range("I122").Value = 1638.4615
$ws.Range("K122").Value = 4915.3845
$ws.Range("M122").Value = -2465.3845
$ws.Range("H126").Value = 2153.1667
$ws.Range("I126").Value = 1806
$ws.Range("K126").Value = 5418
$ws.Range("M126").Value = -2948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 29999
$ws.Range("I45").Value = 29999
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 29999
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = $null
$ws.Range("N45").Value = -29592
$ws.Range("H93").Value = 755.4
$ws.Range("I93").Value = 755.4
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 755.4
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = $null
$ws.Range("N93").Value = 492.6
$ws.Range("H100").Value = 2031.4445
$ws.Range("I100").Value = 1897.1666
$ws.Range("K100").Value = 1897.1666
$ws.Range("M100").Value = -1356.1666
$ws.Range("H132").Value = 34410.16
$ws.Range("I132").Value = 1815.6154
$ws.Range("J132").Value = 57950.668
$ws.Range("K132").Value = 5446.8462
$ws.Range("L132").Value = 173852.004
$ws.Range("M132").Value = -2916.8462
$ws.Range("N132").Value = -178912.004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 15399.667
$ws.Range("J52").Value = 16599.5
$ws.Range("L52").Value = 16599.5
$ws.Range("N52").Value = -17051.5
$ws.Range("H100").Value = 412.36365
$ws.Range("I100").Value = 400.66666
$ws.Range("J100").Value = 426.4
$ws.Range("K100").Value = 801.33332
$ws.Range("L100").Value = 852.8
$ws.Range("M100").Value = -260.33332
$ws.Range("N100").Value = -1934.8
$ws.Range("H132").Value = 6745.5
$ws.Range("I132").Value = 7093
$ws.Range("J132").Value = 6224.25
$ws.Range("K132").Value = 21279
$ws.Range("L132").Value = 18672.75
$ws.Range("M132").Value = -18749
$ws.Range("N132").Value = -23732.75
$ws.Range("H136").Value = 721.6667
$ws.Range("I136").Value = 636.4286
$ws.Range("J136").Value = 1020
$ws.Range("K136").Value = 1909.2858
$ws.Range("L136").Value = 3060
$ws.Range("M136").Value = 640.7142000000001
$ws.Range("N136").Value = -8160
